$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 116, shifting existing rows 116:189 down to 117:190.
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row 116 with its data (same shape as the surrounding rows).
$ws.Range("A116").Value = 4
$ws.Range("B116").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C116").Value = "Los Lagos"
$ws.Range("D116").Value = 44518
$ws.Range("E116").Value = 10
$ws.Range("F116").Value = 100112037
$ws.Range("G116").Value = "Cebollín"
$ws.Range("H116").Value = "Sin especificar"
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 60
$ws.Range("K116").Value = 5000
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = 5000
$ws.Range("N116").Value = "$/paquete 36 unidades"
$ws.Range("O116").Value = "Región Metropolitana"
$ws.Range("P116").Value = 139
$ws.Range("Q116").Value = 36
$ws.Range("R116").Value = "Hortaliza"
